$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (A1:D1) -- bold, thin box border, centered/top-aligned text
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Workspace"
$ws.Range("B1").Value = "Empresa"
$ws.Range("C1").Value = "Data Atualização Dados"
$ws.Range("D1").Value = "Data Verificação"

$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Borders.LineStyle = 1
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160

# Propagate the exact same style to the rest of the header without generating
# extra intermediate cell-style records.
$a1.Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Data rows
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "me"
$ws.Range("B2").Value = "LUPI - Plataforma 2D"
$ws.Range("C2").Value = "24/02/2024"

$ws.Range("A3").Value = "me"
$ws.Range("B3").Value = "Latão Bouticão - Plataforma 2D"
$ws.Range("C3").Value = "24/02/2024"

$ws.Range("A4").Value = "me"
$ws.Range("B4").Value = "MGS - Plataforma 2D"
$ws.Range("C4").Value = "24/02/2024"

# "Data Verificação" column: apply the lowercase variant first, then the
# uppercase one, on the same cell -- this registers both custom number
# formats (164, then 165) while the cell itself ends up carrying only the
# final (165) format, matching the workbook's custom-format history.
$ws.Range("D2").Value = 45348.43799728898
$ws.Range("D2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("D3").Value = 45348.43808098907
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("D4").Value = 45348.43817658117
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
